# Auto-generated edit script applying Tiamat_Profits profit-table updates
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 999
$ws.Range("I40").Value = 999
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 999
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = -824
$ws.Range("N40").ClearContents()

$ws.Range("H64").Value = 18383.705
$ws.Range("I64").Value = 3668.3333
$ws.Range("J64").Value = 128749
$ws.Range("K64").Value = 3668.3333
$ws.Range("L64").Value = 128749
$ws.Range("M64").Value = -3420.3333
$ws.Range("N64").Value = -129245

$ws.Range("H67").Value = 18383.705
$ws.Range("I67").Value = 3668.3333
$ws.Range("J67").Value = 128749
$ws.Range("K67").Value = 3668.3333
$ws.Range("L67").Value = 128749
$ws.Range("M67").Value = -2810.3333
$ws.Range("N67").Value = -130465

$ws.Range("H69").Value = 2453396.2
$ws.Range("I69").Value = 4903460.5
$ws.Range("J69").Value = 3331.6667
$ws.Range("K69").Value = 14710381.5
$ws.Range("L69").Value = 9995.000100000001
$ws.Range("M69").Value = -14709507.5
$ws.Range("N69").Value = -11743.0001

$ws.Range("H72").Value = 2453396.2
$ws.Range("I72").Value = 4903460.5
$ws.Range("J72").Value = 3331.6667
$ws.Range("K72").Value = 44131144.5
$ws.Range("L72").Value = 29985.0003
$ws.Range("M72").Value = -44126776.5
$ws.Range("N72").Value = -38721.0003

$ws.Range("H80").Value = 3498254
$ws.Range("I80").Value = 1368.8667
$ws.Range("J80").Value = 8266734
$ws.Range("K80").Value = 4106.6001
$ws.Range("L80").Value = 24800202
$ws.Range("M80").Value = -3108.6001
$ws.Range("N80").Value = -24802198

$ws.Range("H82").Value = 2686
$ws.Range("I82").Value = 771
$ws.Range("K82").Value = 2313
$ws.Range("M82").Value = -1907

$ws.Range("H83").Value = 3498254
$ws.Range("I83").Value = 1368.8667
$ws.Range("J83").Value = 8266734
$ws.Range("K83").Value = 12319.8003
$ws.Range("L83").Value = 74400606
$ws.Range("M83").Value = -7327.800300000001
$ws.Range("N83").Value = -74410590

$ws.Range("H85").Value = 2686
$ws.Range("I85").Value = 771
$ws.Range("K85").Value = 2313
$ws.Range("M85").Value = -909

$ws.Range("H137").Value = 27988.928
$ws.Range("I137").Value = 36658.82
$ws.Range("J137").Value = 9315.308000000001
$ws.Range("K137").Value = 109976.46
$ws.Range("L137").Value = 27945.924
$ws.Range("M137").Value = -107426.46
$ws.Range("N137").Value = -33045.924

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 171615.14
$ws.Range("I32").Value = 167468.14
$ws.Range("K32").Value = 167468.14
$ws.Range("M32").Value = -167181.14

$ws.Range("H63").Value = 2183.2
$ws.Range("I63").Value = 2336.6667
$ws.Range("J63").Value = 1953
$ws.Range("K63").Value = 2336.6667
$ws.Range("L63").Value = 1953
$ws.Range("M63").Value = -1650.6667
$ws.Range("N63").Value = -3325

$ws.Range("H66").Value = 2183.2
$ws.Range("I66").Value = 2336.6667
$ws.Range("J66").Value = 1953
$ws.Range("K66").Value = 11683.3335
$ws.Range("L66").Value = 9765
$ws.Range("M66").Value = -8251.333500000001
$ws.Range("N66").Value = -16629

$ws.Range("H88").Value = 23939.9
$ws.Range("I88").Value = 10700
$ws.Range("J88").Value = 27249.875
$ws.Range("K88").Value = 10700
$ws.Range("L88").Value = 27249.875
$ws.Range("M88").Value = -10294
$ws.Range("N88").Value = -28061.875

$ws.Range("H91").Value = 23939.9
$ws.Range("I91").Value = 10700
$ws.Range("J91").Value = 27249.875
$ws.Range("K91").Value = 10700
$ws.Range("L91").Value = 27249.875
$ws.Range("M91").Value = -9296
$ws.Range("N91").Value = -30057.875

$ws.Range("H122").Value = 1077.4736
$ws.Range("I122").Value = 978
$ws.Range("J122").Value = 1214.25
$ws.Range("K122").Value = 2934
$ws.Range("L122").Value = 3642.75
$ws.Range("M122").Value = -484
$ws.Range("N122").Value = -8542.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 701590
$ws.Range("I86").Value = 1916.6666
$ws.Range("J86").Value = 1751100
$ws.Range("K86").Value = 1916.6666
$ws.Range("L86").Value = 1751100
$ws.Range("M86").Value = -793.6666
$ws.Range("N86").Value = -1753346

$ws.Range("H89").Value = 701590
$ws.Range("I89").Value = 1916.6666
$ws.Range("J89").Value = 1751100
$ws.Range("K89").Value = 9583.333000000001
$ws.Range("L89").Value = 8755500
$ws.Range("M89").Value = -3967.333000000001
$ws.Range("N89").Value = -8766732

$ws.Range("H105").Value = 1622.5
$ws.Range("I105").Value = 1496.6666
$ws.Range("J105").Value = 2000
$ws.Range("K105").Value = 1496.6666
$ws.Range("L105").Value = 2000
$ws.Range("M105").Value = 250.3334
$ws.Range("N105").Value = -5494

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 27485.3
$ws.Range("I31").Value = 53628.844
$ws.Range("J31").Value = 3831.6191
$ws.Range("K31").Value = 53628.844
$ws.Range("L31").Value = 3831.6191
$ws.Range("M31").Value = -53333.844
$ws.Range("N31").Value = -4421.6191

$ws.Range("H34").Value = 27485.3
$ws.Range("I34").Value = 53628.844
$ws.Range("J34").Value = 3831.6191
$ws.Range("K34").Value = 53628.844
$ws.Range("L34").Value = 3831.6191
$ws.Range("M34").Value = -53426.844
$ws.Range("N34").Value = -4235.6191

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 627.7826
$ws.Range("I113").Value = 484.33334
$ws.Range("J113").Value = 649.3
$ws.Range("K113").Value = 1453.00002
$ws.Range("L113").Value = 1947.9
$ws.Range("M113").Value = 716.9999800000001
$ws.Range("N113").Value = -6287.9

$ws.Range("H128").Value = 430000
$ws.Range("I128").Value = 430000
$ws.Range("K128").Value = 1290000
$ws.Range("M128").Value = -1285020

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4241.5
$ws.Range("I70").Value = 4179.8
$ws.Range("J70").Value = 4550
$ws.Range("K70").Value = 4179.8
$ws.Range("L70").Value = 4550
$ws.Range("M70").Value = -3909.8
$ws.Range("N70").Value = -5090

$ws.Range("H73").Value = 4241.5
$ws.Range("I73").Value = 4179.8
$ws.Range("J73").Value = 4550
$ws.Range("K73").Value = 4179.8
$ws.Range("L73").Value = 4550
$ws.Range("M73").Value = -3243.8
$ws.Range("N73").Value = -6422

$ws.Range("H97").Value = 1261.6666
$ws.Range("I97").Value = 1159.5555
$ws.Range("K97").Value = 1159.5555
$ws.Range("M97").Value = -663.5554999999999

$ws.Range("H102").Value = 35690
$ws.Range("I102").Value = 15058.25
$ws.Range("K102").Value = 15058.25
$ws.Range("M102").Value = -13436.25

$ws.Range("H126").Value = 1883.1765
$ws.Range("I126").Value = 1400
$ws.Range("J126").Value = 2221.4
$ws.Range("K126").Value = 4200
$ws.Range("L126").Value = 6664.200000000001
$ws.Range("M126").Value = -1730
$ws.Range("N126").Value = -11604.2

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2953.25
$ws.Range("I46").Value = 410.5
$ws.Range("J46").Value = 5496
$ws.Range("K46").Value = 410.5
$ws.Range("L46").Value = 5496
$ws.Range("M46").Value = -222.5
$ws.Range("N46").Value = -5872

$ws.Range("H81").Value = 30000
$ws.Range("I81").Value = 0
$ws.Range("J81").Value = 30000
$ws.Range("K81").Value = 0
$ws.Range("L81").Value = 30000
$ws.Range("M81").ClearContents()
$ws.Range("N81").Value = -31996

$ws.Range("H82").Value = 2185.3572
$ws.Range("I82").Value = 2200
$ws.Range("J82").Value = 2179.5
$ws.Range("K82").Value = 2200
$ws.Range("L82").Value = 2179.5
$ws.Range("M82").Value = -1839
$ws.Range("N82").Value = -2901.5

$ws.Range("H84").Value = 30000
$ws.Range("I84").Value = 0
$ws.Range("J84").Value = 30000
$ws.Range("K84").Value = 0
$ws.Range("L84").Value = 90000
$ws.Range("M84").ClearContents()
$ws.Range("N84").Value = -99984

$ws.Range("H85").Value = 2185.3572
$ws.Range("I85").Value = 2200
$ws.Range("J85").Value = 2179.5
$ws.Range("K85").Value = 2200
$ws.Range("L85").Value = 2179.5
$ws.Range("M85").Value = -952
$ws.Range("N85").Value = -4675.5

$ws.Range("H87").Value = 30000
$ws.Range("J87").Value = 30000
$ws.Range("L87").Value = 30000
$ws.Range("N87").Value = -32246

$ws.Range("H90").Value = 30000
$ws.Range("J90").Value = 30000
$ws.Range("L90").Value = 90000
$ws.Range("N90").Value = -101232

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H76").Value = 15250
$ws.Range("I76").Value = 10000
$ws.Range("K76").Value = 10000
$ws.Range("M76").Value = -9685

$ws.Range("H79").Value = 15250
$ws.Range("I79").Value = 10000
$ws.Range("K79").Value = 10000
$ws.Range("M79").Value = -8908

$ws.Range("H80").Value = 45250
$ws.Range("J80").Value = 45250
$ws.Range("L80").Value = 45250
$ws.Range("N80").Value = -47246

$ws.Range("H83").Value = 45250
$ws.Range("J83").Value = 45250
$ws.Range("L83").Value = 135750
$ws.Range("N83").Value = -145734
